# Auto-generated Excel COM-interop script
# Applies numeric cell updates (profit recalculation) across all 8 item-sheets
# as captured by the upstream diff for Sheets/Malboro_Profits.xlsx

$wb = $excel.ActiveWorkbook


# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 348.75
$ws.Range("J2").Value = 367.85715
$ws.Range("L2").Value = 367.85715
$ws.Range("N2").Value = -593.85715
$ws.Range("H43").Value = 10400
$ws.Range("I43").Value = 10133.333
$ws.Range("K43").Value = 10133.333
$ws.Range("M43").Value = -10064.333
$ws.Range("H64").Value = 6749.5
$ws.Range("J64").Value = 6749.5
$ws.Range("L64").Value = 6749.5
$ws.Range("N64").Value = -7245.5
$ws.Range("H67").Value = 6749.5
$ws.Range("J67").Value = 6749.5
$ws.Range("L67").Value = 6749.5
$ws.Range("N67").Value = -8465.5
$ws.Range("H86").Value = 18736.2
$ws.Range("I86").Value = 15185.737
$ws.Range("J86").Value = 29979.334
$ws.Range("K86").Value = 15185.737
$ws.Range("L86").Value = 29979.334
$ws.Range("M86").Value = -14062.737
$ws.Range("N86").Value = -32225.334
$ws.Range("H88").Value = 1926.1111
$ws.Range("J88").Value = 1440.1818
$ws.Range("L88").Value = 1440.1818
$ws.Range("N88").Value = -2252.1818
$ws.Range("H89").Value = 18736.2
$ws.Range("I89").Value = 15185.737
$ws.Range("J89").Value = 29979.334
$ws.Range("K89").Value = 75928.685
$ws.Range("L89").Value = 149896.67
$ws.Range("M89").Value = -70312.685
$ws.Range("N89").Value = -161128.67
$ws.Range("H91").Value = 1926.1111
$ws.Range("J91").Value = 1440.1818
$ws.Range("L91").Value = 1440.1818
$ws.Range("N91").Value = -4248.1818
$ws.Range("H98").Value = 2159.1035
$ws.Range("I98").Value = 2245.36
$ws.Range("K98").Value = 2245.36
$ws.Range("M98").Value = -747.3600000000001
$ws.Range("H122").Value = 2159.1035
$ws.Range("I122").Value = 2245.36
$ws.Range("K122").Value = 6736.08
$ws.Range("M122").Value = -4286.08
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 9824.457
$ws.Range("I132").Value = 7353.0303
$ws.Range("J132").Value = 50603
$ws.Range("K132").Value = 22059.0909
$ws.Range("L132").Value = 151809
$ws.Range("M132").Value = -19529.0909
$ws.Range("N132").Value = -156869
$ws.Range("H137").Value = 10435.263
$ws.Range("I137").Value = 2283.7036
$ws.Range("J137").Value = 30443.637
$ws.Range("K137").Value = 6851.110799999999
$ws.Range("L137").Value = 91330.91099999999
$ws.Range("M137").Value = -4301.110799999999
$ws.Range("N137").Value = -96430.91099999999

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 3307.43
$ws.Range("I32").Value = 1886.3478
$ws.Range("K32").Value = 1886.3478
$ws.Range("M32").Value = -1599.3478
$ws.Range("H61").Value = 16507.44
$ws.Range("J61").Value = 19634.867
$ws.Range("L61").Value = 19634.867
$ws.Range("N61").Value = -20058.867
$ws.Range("H132").Value = 2867335
$ws.Range("I132").Value = 4180.2
$ws.Range("K132").Value = 12540.6
$ws.Range("M132").Value = -10010.6
$ws.Range("H136").Value = 16507.44
$ws.Range("J136").Value = 19634.867
$ws.Range("L136").Value = 58904.601
$ws.Range("N136").Value = -64004.601

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 1808.4
$ws.Range("I86").Value = 1904.2354
$ws.Range("J86").Value = 1265.3334
$ws.Range("K86").Value = 1904.2354
$ws.Range("L86").Value = 1265.3334
$ws.Range("M86").Value = -781.2354
$ws.Range("N86").Value = -3511.3334
$ws.Range("H89").Value = 1808.4
$ws.Range("I89").Value = 1904.2354
$ws.Range("J89").Value = 1265.3334
$ws.Range("K89").Value = 9521.177
$ws.Range("L89").Value = 6326.666999999999
$ws.Range("M89").Value = -3905.177
$ws.Range("N89").Value = -17558.667
$ws.Range("H94").Value = 4758.5
$ws.Range("I94").Value = 4124.974
$ws.Range("K94").Value = 4124.974
$ws.Range("M94").Value = -3673.974

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 4995.9287
$ws.Range("I16").Value = 786.4
$ws.Range("J16").Value = 7334.5557
$ws.Range("K16").Value = 786.4
$ws.Range("L16").Value = 7334.5557
$ws.Range("M16").Value = -499.4
$ws.Range("N16").Value = -7908.5557
$ws.Range("H28").Value = 18880.666
$ws.Range("J28").Value = 18880.666
$ws.Range("L28").Value = 18880.666
$ws.Range("N28").Value = -19370.666
$ws.Range("H31").Value = 48428.6
$ws.Range("I31").Value = 19065
$ws.Range("J31").Value = 68004.336
$ws.Range("K31").Value = 19065
$ws.Range("L31").Value = 68004.336
$ws.Range("M31").Value = -18770
$ws.Range("N31").Value = -68594.336
$ws.Range("H34").Value = 48428.6
$ws.Range("I34").Value = 19065
$ws.Range("J34").Value = 68004.336
$ws.Range("K34").Value = 19065
$ws.Range("L34").Value = 68004.336
$ws.Range("M34").Value = -18863
$ws.Range("N34").Value = -68408.336
$ws.Range("H97").Value = 16800
$ws.Range("J97").Value = 16800
$ws.Range("L97").Value = 16800
$ws.Range("N97").Value = -18782
$ws.Range("H113").Value = 4995.9287
$ws.Range("I113").Value = 786.4
$ws.Range("J113").Value = 7334.5557
$ws.Range("K113").Value = 786.4
$ws.Range("L113").Value = 7334.5557
$ws.Range("M113").Value = 1383.6
$ws.Range("N113").Value = -11674.5557
$ws.Range("H131").Value = 59990
$ws.Range("J131").Value = 59990
$ws.Range("L131").Value = 59990
$ws.Range("N131").Value = -70070
$ws.Range("H132").Value = 8891.237999999999
$ws.Range("I132").Value = 2441.3845
$ws.Range("K132").Value = 7324.1535
$ws.Range("M132").Value = -4794.1535
$ws.Range("H141").Value = 204074.36
$ws.Range("J141").Value = 204074.36
$ws.Range("L141").Value = 204074.36
$ws.Range("N141").Value = -214434.36

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H99").Value = 2725.5
$ws.Range("I99").Value = 733
$ws.Range("K99").Value = 2199
$ws.Range("M99").Value = 47
$ws.Range("H121").Value = 612
$ws.Range("I121").Value = 366.42856
$ws.Range("J121").Value = 803
$ws.Range("K121").Value = 1099.28568
$ws.Range("L121").Value = 2409
$ws.Range("M121").Value = 210.71432
$ws.Range("N121").Value = -5029
$ws.Range("H131").Value = 1498.07
$ws.Range("J131").Value = 1498.07
$ws.Range("L131").Value = 4494.21
$ws.Range("N131").Value = -14574.21

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 5215.7896
$ws.Range("I102").Value = 5594.2354
$ws.Range("K102").Value = 5594.2354
$ws.Range("M102").Value = -3972.2354
$ws.Range("H132").Value = 10631.286
$ws.Range("I132").Value = 10722.85
$ws.Range("K132").Value = 32168.55
$ws.Range("M132").Value = -29638.55

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H22").Value = 10590.913
$ws.Range("I22").Value = 9784.076999999999
$ws.Range("K22").Value = 9784.076999999999
$ws.Range("M22").Value = -9489.076999999999
$ws.Range("H27").Value = 10590.913
$ws.Range("I27").Value = 9784.076999999999
$ws.Range("K27").Value = 9784.076999999999
$ws.Range("M27").Value = -9677.076999999999
$ws.Range("H132").Value = 4472000
$ws.Range("I132").Value = 2999
$ws.Range("J132").Value = 5030625
$ws.Range("K132").Value = 8997
$ws.Range("L132").Value = 15091875
$ws.Range("M132").Value = -6467
$ws.Range("N132").Value = -15096935

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H38").Value = 19874.5
$ws.Range("I38").Value = 19874.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 19874.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -19401.5
$ws.Range("N38").ClearContents()
$ws.Range("H126").Value = 5491.609
$ws.Range("I126").Value = 4233.4375
$ws.Range("J126").Value = 8367.429
$ws.Range("K126").Value = 12700.3125
$ws.Range("L126").Value = 25102.287
$ws.Range("M126").Value = -10230.3125
$ws.Range("N126").Value = -30042.287
$ws.Range("H132").Value = 12513.762
$ws.Range("I132").Value = 2857.3333
$ws.Range("K132").Value = 8571.999899999999
$ws.Range("M132").Value = -6041.999899999999
$ws.Range("H136").Value = 11326.885
$ws.Range("J136").Value = 17986.732
$ws.Range("L136").Value = 53960.196
$ws.Range("N136").Value = -59060.196
$ws.Range("H140").Value = 129156.664
$ws.Range("J140").Value = 129156.664
$ws.Range("L140").Value = 129156.664
$ws.Range("N140").Value = -139516.664
